# Update CalSim3 init file: scenario 25 -> scenario 27 related cell references
# (the commit message says "updated init files for scenarios 25 and 27")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# D5:D11 hold "Lower Right Cell" references that move from row 26 to row 28
$ws.Range("D5").Value = "A28"
$ws.Range("D6").Value = "B28"
$ws.Range("D7").Value = "C28"
$ws.Range("D8").Value = "G28"
$ws.Range("D9").Value = "H28"
$ws.Range("D10").Value = "I28"
$ws.Range("D11").Value = "J28"

# Update the active selection to match the edited file (D11 now selected instead of D12)
$ws.Range("D11").Select()
